$d = $word.ActiveDocument

# Locate the "Requisitos" section's requirement paragraph, then remove the
# trailing empty paragraph, the "Ver no Jupiter..." paragraph, and the
# "(c) 2020..." footer paragraph that follow it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "LOM3204: Desenho Técnico e Projeto Assistido por Computador (Requisito)") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $startPara = $d.Paragraphs.Item($target + 1)
    $endPara = $d.Paragraphs.Item($target + 3)
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
